# Update GitHub org URL
# (1) Footer "last updated" date placeholders (slide master + all layouts):
#     2019-08-28 -> 2019-09-13
# (2) GitHub org URL text runs: "SSC-CIO-DigitalTaskForce" -> "dtf-ein"
#     on the [GET IN TOUCH] slide and the [JOIN IN] slide.
# (3) ENGAGE tag bookkeeping Id refreshed for the new save.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer date placeholders on the slide master and every layout.
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "2019-09-13"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) GitHub org URL text on slide 9 ("[GET IN TOUCH]") and slide 11
#    ("[JOIN IN]"). Replace "SSC-CIO-" + "DigitalTaskForce" with
#    "dtf-" + "ein" while keeping existing run formatting intact.
# ---------------------------------------------------------------------
function Update-GithubOrg($tr) {
    $full = $tr.Text
    $idx = $full.IndexOf("github.com/SSC-CIO-")
    if ($idx -ge 0) {
        $start = $idx + 1
        $len = [string]"github.com/SSC-CIO-".Length
        $tr.Characters($start, $len).Text = "github.com/dtf-"
    }

    $full = $tr.Text
    $idx = $full.IndexOf("DigitalTaskForce")
    if ($idx -ge 0) {
        $start = $idx + 1
        $len = [string]"DigitalTaskForce".Length
        $tr.Characters($start, $len).Text = "ein"
    }
}

# Slide 9: Group 13 > TextBox 4 ("GitHub:	github.com/SSC-CIO-DigitalTaskForce...")
$slide9 = $p.Slides.Item(9)
$group13 = $slide9.Shapes.Item(3)
$textBox4 = $group13.GroupItems.Item(1)
Update-GithubOrg $textBox4.TextFrame.TextRange

# Slide 11: Rectangle 49 ("[JOIN IN] ... - github.com/SSC-CIO-DigitalTaskForce...")
$slide11 = $p.Slides.Item(11)
$rect49 = $slide11.Shapes.Item(18)
Update-GithubOrg $rect49.TextFrame.TextRange

# ---------------------------------------------------------------------
# 3) ENGAGE tag: refresh the opaque bookkeeping Id embedded in the JSON
#    blob stored on the presentation tag.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $p.Tags.Count; $i++) {
    if ($p.Tags.Name($i) -eq "ENGAGE") {
        $val = $p.Tags.Value($i)
        $newVal = $val.Replace("5d64385d30303426cc9a71a7", "5d7bd6ac4331434bc0c5002d")
        $p.Tags.Add("ENGAGE", $newVal)
    }
}
